$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-19 Sunday" "2025-01-20 Monday"

Replace-Text "184×3=" "365×8="
Replace-Text "312×7=" "752×4="
Replace-Text "596×7=" "846×7="
Replace-Text "447×5=" "141×9="
Replace-Text "794×6=" "608×7="
Replace-Text "593×5=" "852×8="
Replace-Text "316×4=" "569×6="
Replace-Text "583×7=" "625×4="
Replace-Text "850×2=" "151×7="
Replace-Text "881×2=" "123×5="
Replace-Text "280×6=" "705×4="
Replace-Text "525×6=" "786×6="
Replace-Text "661×6=" "441×8="
Replace-Text "683×6=" "923×4="
Replace-Text "851×9=" "506×4="
Replace-Text "224×6=" "448×5="
Replace-Text "347×5=" "504×8="
Replace-Text "369×2=" "426×2="
Replace-Text "726×3=" "474×3="
Replace-Text "869×5=" "828×7="
Replace-Text "107×4=" "536×8="
Replace-Text "577×6=" "911×2="
Replace-Text "978×5=" "929×7="
Replace-Text "385×2=" "133×5="
Replace-Text "699×3=" "546×8="
